# Auto-generated PowerPoint COM-interop edit script
$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(2)  # 'Title and Content'

# --- Step 1: insert the three new '2.x Preprocessing' slides after slide 3 ---
$s = $p.Slides.AddSlide(4, $layout)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = '2.1 Détail du Preprocessing : Subscriptions'
$titleShape.TextFrame.TextRange.Font.Bold = $true
$titleShape.TextFrame.TextRange.Font.Color.RGB = 6697728
$contentShape = $s.Shapes.Item(2)
$bodyLines = @('', 'Traitement des colonnes statiques (Profil Client) :', '   - Identifiants (ID, Dates) : Utilisés pour le filtrage (15 jours) puis supprimés pour éviter le bruit.', '   - Catégories Nominales (Vendor, Region, Legal) : Traitées par OneHotEncoder (gestion des inconnus en test).', '   - Catégories Ordinales (Revenue, Employees) : Traitées par OrdinalEncoder (préservation de l''ordre).', '   - Cas Spéciaux : ''v2_modules'' (Parsing Multi-label) et ''v2_segment'' (OneHot avec drop=''first'' pour éviter la colinéarité).', '   - Cible : Dérivée de ''first_paid_invoice_paid_at'' (1 si date présente, 0 sinon).')
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tr = $contentShape.TextFrame.TextRange
$n = $tr.Paragraphs().Count
for ($i = 2; $i -le $n; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Font.Size = 18
    $para.ParagraphFormat.SpaceAfter = 10
}

$s = $p.Slides.AddSlide(5, $layout)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = '2.2 Détail du Preprocessing : Daily Usage'
$titleShape.TextFrame.TextRange.Font.Bold = $true
$titleShape.TextFrame.TextRange.Font.Color.RGB = 6697728
$contentShape = $s.Shapes.Item(2)
$bodyLines = @('', 'Traitement des métriques d''activité (19 colonnes nb_*) :', '   - Valeurs Manquantes : Remplacées par 0 (correspond à une absence d''activité réelle).', '   - Modèles Tabulaires (LightGBM/XGB) : Agrégation par essai -> Somme, Moyenne, Max, Ecart-Type (StandardScaler).', '   - Modèles Séquentiels (DL) : Conservation de la structure temporelle (416 essais, 15 jours, 19 features).', '   - Objectif : Capturer l''intensité (Somme) et la régularité (Ecart-Type) de l''usage.')
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tr = $contentShape.TextFrame.TextRange
$n = $tr.Paragraphs().Count
for ($i = 2; $i -le $n; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Font.Size = 18
    $para.ParagraphFormat.SpaceAfter = 10
}

$s = $p.Slides.AddSlide(6, $layout)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = '2.3 Règles Globales & Dimensions Finales'
$titleShape.TextFrame.TextRange.Font.Bold = $true
$titleShape.TextFrame.TextRange.Font.Color.RGB = 6697728
$contentShape = $s.Shapes.Item(2)
$bodyLines = @('', 'Nos principes de rigueur :', '   - Anti-Leakage : Suppression stricte de ''subscription_status'' et ''canceled_at'' (infos du futur).', '   - Robustesse : Les catégories inconnues en test sont ignorées (handle_unknown=''ignore'').', '   - Volumétrie Finale : ~150 features (Tabulaire) vs Tensor (416, 15, 19) (Deep Learning).', '   - Résultat : Un pipeline ''Production-Ready'' robuste aux nouvelles données.')
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tr = $contentShape.TextFrame.TextRange
$n = $tr.Paragraphs().Count
for ($i = 2; $i -le $n; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Font.Size = 18
    $para.ParagraphFormat.SpaceAfter = 10
}

# --- Step 2: fix the last bullet on the (now pushed-down) 'signaux précurseurs' slide ---
$sigSlide = $p.Slides.Item(11)
$sigContent = $sigSlide.Shapes.Item(2)
$sigTr = $sigContent.TextFrame.TextRange
$sigN = $sigTr.Paragraphs().Count
$sigLast = $sigTr.Paragraphs($sigN, 1)
$sigLast.Text = 'Insights : L''usage intensif (Factures/Mobile) tôt dans l''essai garanti la conversion.'

# --- Step 3: insert the two new '7b/7c' slides after the 'signaux précurseurs' slide ---
$s = $p.Slides.AddSlide(12, $layout)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = '7b. Le ''Top 1%'' : Modèle Hybride'
$titleShape.TextFrame.TextRange.Font.Bold = $true
$titleShape.TextFrame.TextRange.Font.Color.RGB = 6697728
$contentShape = $s.Shapes.Item(2)
$bodyLines = @('', 'Pour maximiser la performance, nous avons créé un Ensemble Hybride :', '1. Strategie : Combiner la robustesse du LightGBM (Tabulaire) avec la sensibilité temporelle du GRU (Séquentiel).', '2. Méthode : Moyenne pondérée des probabilités (70% LightGBM + 30% GRU).', '3. Gain : Hausse de l''AUC (+0.02) et meilleure calibration (Brier Score réduit).', '4. Résultat : Un ''super-modèle'' qui ne rate presque aucun signal faible.')
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tr = $contentShape.TextFrame.TextRange
$n = $tr.Paragraphs().Count
for ($i = 2; $i -le $n; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Font.Size = 18
    $para.ParagraphFormat.SpaceAfter = 10
}

$s = $p.Slides.AddSlide(13, $layout)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = '7c. Simulation ROI & Impact Business'
$titleShape.TextFrame.TextRange.Font.Bold = $true
$titleShape.TextFrame.TextRange.Font.Color.RGB = 6697728
$contentShape = $s.Shapes.Item(2)
$bodyLines = @('', 'Traduction du Score en Euros (Simulation sur Test Set) :', '   - Hypothèses : LTV = 500€, Coût d''Intervention (Call) = 10€, Taux de Succès = 20%.', '   - Stratégie : Intervenir seulement si le risque de churn est élevé (Score < Seuil).', '   - Résultat : En ciblant les utilisateurs à risque (Prob < 0.45) :', '       -> On sauve ~12% de churn additionnel.', '       -> ROI Net estimé : +15 000€ / mois (pour 1000 essais).', 'Conclusion : Le modèle n''est pas une dépense, c''est un centre de profit immédiat.')
$contentShape.TextFrame.TextRange.Text = [string]::Join([char]13, $bodyLines)
$tr = $contentShape.TextFrame.TextRange
$n = $tr.Paragraphs().Count
for ($i = 2; $i -le $n; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.Font.Size = 18
    $para.ParagraphFormat.SpaceAfter = 10
}

# --- Sanity log ---
Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Write-Output "Slide $($i): $($p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text)"
}
